$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the naive-forecast cells C2 and C3 (bug: these had stray leaked values)
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Corrected forecast values (naive component forecaster bug fix)
$ws.Range("E5").Value = 2.743551941645195
$ws.Range("E6").Value = 1.794132456841191
$ws.Range("C7").Value = 1.239479831392831
$ws.Range("C8").Value = 0.2379616621360992
$ws.Range("E8").Value = 1.348985046565399
$ws.Range("E9").Value = 0.7739869831244084
$ws.Range("C10").Value = 1.470039379455734
$ws.Range("E10").Value = 1.319057785023636
$ws.Range("C11").Value = 1.638797242243228
$ws.Range("E11").Value = 1.236938064849946
$ws.Range("E12").Value = 2.320541194291859
$ws.Range("E13").Value = 0.6374750548025609
$ws.Range("E14").Value = 0.5721424576528022
$ws.Range("C16").Value = 1.099928004397577
$ws.Range("E16").Value = 1.661735028110978
$ws.Range("C17").Value = 2.310042359896247
$ws.Range("E17").Value = 1.255295035968373
$ws.Range("C19").Value = -0.3101476031197037
